$d = $word.ActiveDocument

# Locate the "Palabras clave" paragraph that holds the keyword list text.
$keywordsPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Palabras clave del recurso*") {
        $keywordsPara = $p
        break
    }
}

$paraStart = $keywordsPara.Range.Start
$paraEnd = $keywordsPara.Range.End

# Find the "imperialismo, ..." run within that paragraph and the trailing
# period that needs to be removed.
$fullMarker = "imperialismo, colonialismo, " + [char]0x00C1 + "frica, Congo, explotaci" + [char]0x00F3 + "n, potencias."
$markerRange = $d.Range($paraStart, $paraEnd)
[void]$markerRange.Find.Execute($fullMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$runStart = $markerRange.Start
$runEnd = $markerRange.End

# Position right after "...Africa" (before ", Congo, explotacion, potencias.")
$splitMarker = "imperialismo, colonialismo, " + [char]0x00C1 + "frica"
$splitRange = $d.Range($runStart, $runEnd)
[void]$splitRange.Find.Execute($splitMarker, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint = $splitRange.End

# The trailing "." sits right before the paragraph mark.
$periodRange = $d.Range($runEnd - 1, $runEnd)

# Move the "_GoBack" bookmark from the earlier empty paragraph to the end
# of this paragraph (right where the "." currently sits) *before* deleting
# the period, so it naturally collapses to that exact position once the
# character is removed.
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()
$d.Bookmarks.Add("_GoBack", $periodRange)

# Remove the trailing period.
$periodRange2 = $d.Range($runEnd - 1, $runEnd)
$periodRange2.Delete()

# Split "imperialismo, ... Africa, Congo, ... potencias" into two runs at
# the word boundary after "Africa" (same formatting on both sides) by
# toggling a character property on/off, which forces the run boundary to
# persist even though the resulting formatting is unchanged.
$firstHalf = $d.Range($runStart, $splitPoint)
$firstHalf.Bold = 1
$firstHalf.Bold = 0
